$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new question rows (27 and 28), mirroring the style/pattern of
# the existing rows (row 26 in particular).
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "Print Nodes in Top View of Binary Tree"
$ws.Range("D27").Value = "Tree"
$ws.Range("E27").Value = "medium"
$ws.Range("F27").Value = "GeeksForGeeks"

$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "Remove nodes on root to leaf paths of length < K"
$ws.Range("D28").Value = "Tree"
$ws.Range("E28").Value = "medium"
$ws.Range("F28").Value = "GeeksForGeeks"

# Match styles used by the rest of the table: column A/D/E/F centered (style 1),
# column B left aligned (style 3).
$ws.Range("A27:A28").HorizontalAlignment = -4108
$ws.Range("D27:F28").HorizontalAlignment = -4108
$ws.Range("B27:B28").HorizontalAlignment = -4131

# Update the selection/view to match the edited state.
$ws.Application.ActiveWindow.Zoom = 115
$ws.Range("C28").Select()
